# "add process cadastro produto"
#
# The product-registration process needs the unit-of-measure ("UND", column C)
# labels to carry a couple of trailing spaces so they line up with the rest of
# the catalogue's category labels (which already end in "  "). Normalize the
# six affected unit labels wherever they appear in the Data sheet, then leave
# the selection where the editor last left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Exact-value map: only touch cells whose current value is *exactly* one of
# these labels (column C only has these 7 possible unit values), so long
# product descriptions that merely contain "KIT" or "PAR" as a substring are
# left untouched.
$unitRenames = @{
    "Caixa"   = "Caixa  "
    "KIT"     = "KIT  "
    "Metro"   = "Metro  "
    "PACOTE"  = "PACOTE  "
    "PAR"     = "PAR  "
    "Unidade" = "Unidade  "
}

$lastRow = $ws.UsedRange.Rows.Count
$colC = 3

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $colC)
    $current = $cell.Value2
    if ($unitRenames.ContainsKey($current)) {
        $cell.Value2 = $unitRenames[$current]
    }
}

# Restore the editor's last selection on the sheet.
$ws.Range("J17").Select()
